# Applies the "traded" update: fills in the PriceChange/UpDown columns for
# the most recent existing row (12) and appends the new day's row (13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Finish row 12 (PriceChange / UpDown were missing) ---
$ws.Cells.Item(12, 24).Value = -0.16000300000000323   # X12 PriceChange
$ws.Cells.Item(12, 25).Value = "Down"                  # Y12 UpDown

# --- Append new row 13 with the latest trading data ---
$ws.Cells.Item(13, 1).Value = 42654.894490740742       # A13 Date
$ws.Cells.Item(13, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(13, 2).Value = 9                         # B13 ScoreFinal
$ws.Cells.Item(13, 3).Value = "Buy"                     # C13 Verdict
$ws.Cells.Item(13, 4).Value = 4                         # D13 totalSentiment
$ws.Cells.Item(13, 5).Value = 21700                     # E13 wordCount
$ws.Cells.Item(13, 6).Value = 2567                      # F13 sentenceCount
$ws.Cells.Item(13, 7).Value = 54                        # G13 posWordPercentage
$ws.Cells.Item(13, 8).Value = 44                        # H13 negWordPercentage
$ws.Cells.Item(13, 9).Value = 73                        # I13 posPhrasePercentage
$ws.Cells.Item(13, 10).Value = 26                       # J13 negPhrasePercentage
$ws.Cells.Item(13, 11).Value = 12049                    # K13 ElapsedMs
$ws.Cells.Item(13, 12).Value = 270                      # L13 posWordCount
$ws.Cells.Item(13, 13).Value = 222                      # M13 negWordCount
$ws.Cells.Item(13, 14).Value = 70                       # N13 positivePhraseCount
$ws.Cells.Item(13, 15).Value = 25                       # O13 negativePhraseCount
$ws.Cells.Item(13, 16).Value = "Bag"                    # P13 Method
$ws.Cells.Item(13, 17).Value = 28.689659976213832       # Q13 RSI
$ws.Cells.Item(13, 18).Value = 0.84                     # R13 PEG

$ws.Cells.Item(13, 19).Value = -0.0125                  # S13 200Moving%
$ws.Cells.Item(13, 19).NumberFormat = "0.00%"

$ws.Cells.Item(13, 20).Value = -0.0261                  # T13 50Moving%
$ws.Cells.Item(13, 20).NumberFormat = "0.00%"

$ws.Cells.Item(13, 21).Value = 14.56                    # U13 PriceBook
$ws.Cells.Item(13, 22).Value = "N/A"                    # V13 Dividend
$ws.Cells.Item(13, 23).Value = 0                        # W13 Bollinger
